########################################################################
# this is first commit
#
# Adds a small header row (s.no / date / team) to Sheet1 in D2:F2,
# draws a thin box border around D2:G2 (G2 stays empty), and leaves the
# selection on E4, matching the authored workbook change.
########################################################################

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header values (become shared strings s.no / date / team)
$ws.Range("D2").Value = "s.no"
$ws.Range("E2").Value = "date"
$ws.Range("F2").Value = "team"

# Thin box border around D2:G2 (G2 is left blank but keeps the style)
$rng = $ws.Range("D2:G2")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Leave the cursor/selection on E4, as in the authored workbook
$sel = $ws.Range("E4").Select()
